# Generate Report for Handoff
# Update the "Latest Handoff Datetime" column for the rows that were just
# handed off (708e369f for zh-cn, 556b5bf8 for de-de).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 5 -> 708e369f...zh-cn.xlf, column D = Latest Handoff Datetime
$wsZhCn.Range("D5").Value = "2016-03-04 01:18:29"

# de-de sheet: row 4 -> 556b5bf8...de-de.xlf, column D = Latest Handoff Datetime
$wsDeDe.Range("D4").Value = "2016-03-04 01:18:43"
